$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update description of row 47 (Sarvan Amel's first ERD task):
# "Creating ERD" -> "Creating first version of ERD"
$ws.Range("B47").Value = "Creating first version of ERD"

# Update the hours formula/value for row 54 ("Adding new Tables to ERD"):
# 57/60 (0.95) -> 86/60 (1.4333333333333333)
$ws.Range("C54").Formula = "=86/60"

# Recalculate dependent formulas (B62 totals, B63 grand total)
$excel.CalculateFullRebuild()
